{"js": "// Replace the outdated description of the loan/reservation management\n// form with the updated text describing the new form-based workflow.\nconst oldText =\n  \"posudbama i rezervacijama pristupa s po\u010detne stranice. U pregledu posudbi \" +\n  \"mo\u017ee pretra\u017eivati posudbe po razli\u010ditim kriterijima i ako je mogu\u0107e \" +\n  \"produljiti posudbu. U upravljanju  posudbama preko skenera mo\u017ee \" +\n  \"skenirati primjerak i karticu korisnika i tako vratiti posuditi ili \" +\n  \"rezervirati primjerak.\";\n\nconst newText =\n  \"posudbama i rezervacijama pristupa s po\u010detne stranice. Na formi posudbe \" +\n  \"ima pregled posudbi koje mo\u017ee pretra\u017eivati i produljiti. Na formi \" +\n  \"Upravljanje posudbama na po\u010detku skenira korisnikov QR kod i onda \" +\n  \"skenira bar kod primjerka. Ako su uneseni svi podaci, primjerak se mo\u017ee \" +\n  \"vratiti/posuditi/rezervirati.\";\n\nconst results = context.document.body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\nresults.items[0].insertText(newText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the outdated description of the loan/reservation management\n# form with the updated text describing the new form-based workflow.\n$d = $word.ActiveDocument\n\n$oldText = \"posudbama i rezervacijama pristupa s po\u010detne stranice. U pregledu posudbi mo\u017ee pretra\u017eivati posudbe po razli\u010ditim kriterijima i ako je mogu\u0107e produljiti posudbu. U upravljanju  posudbama preko skenera mo\u017ee skenirati primjerak i karticu korisnika i tako vratiti posuditi ili rezervirati primjerak.\"\n$newText = \"posudbama i rezervacijama pristupa s po\u010detne stranice. Na formi posudbe ima pregled posudbi koje mo\u017ee pretra\u017eivati i produljiti. Na formi Upravljanje posudbama na po\u010detku skenira korisnikov QR kod i onda skenira bar kod primjerka. Ako su uneseni svi podaci, primjerak se mo\u017ee vratiti/posuditi/rezervirati.\"\n\n$targetRange = $d.Content\n$find = $targetRange.Find\n$find.Text = $oldText\n$found = $find.Execute($oldText)\nif (-not $found) {\n    throw \"Target text not found\"\n}\n\n# Using FormattedText (instead of plain Range.Text / Find-replace) keeps the\n# untouched preceding run (\"Moderator upravljanju \") intact instead of being\n# coalesced into the edited run.\n$ft = $targetRange.FormattedText\n$ft.Text = $newText\n\n$freshRange = $d.Range($ft.Start, $ft.End)\n$freshRange.FormattedText = $ft\n"}
